$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto-price refresh. The Price/Volume columns (D:E) are stored as
# plain text (e.g. "1.000", "21.80", "  +0.65%  ") rather than numbers, so force a
# text number-format before writing; this stops Excel from auto-coercing values
# that merely look numeric (e.g. "1.000" -> 1). The format gets reset back to the
# workbook default ("Normal" style) afterwards so cell styling is left untouched.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.231.12"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "1.860.39"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "236.64"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "0.4674"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.2865"
$ws.Range("E8").Value = "  +0.92%  "

$ws.Range("D9").Value = "0.06542"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").Value = "21.80"
$ws.Range("E10").Value = "  +4.98%  "

$ws.Range("D11").Value = "0.07921"
$ws.Range("E11").Value = "  +0.60%  "

$ws.Range("D12").Value = "97.58"
$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").Value = "1.866.33"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").Value = "5.175"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").Value = "0.6796"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").Value = "267.80"
$ws.Range("E16").Value = "  -4.68%  "

$ws.Range("D17").Value = "30.218.07"

$ws.Range("D18").Value = "13.84"
$ws.Range("E18").Value = "  +9.36%  "

$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "0.000007394"
$ws.Range("E20").Value = "  +1.82%  "

$ws.Range("D21").Value = "2.111.23"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "5.323"
$ws.Range("E22").Value = "  -3.17%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").Value = "6.196"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").Value = "167.19"
$ws.Range("E25").Value = "  +1.25%  "

$ws.Range("D26").Value = "9.221"
$ws.Range("E26").Value = "  -0.88%  "

$ws.Range("D27").Value = "18.89"
$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("D28").Value = "1.972"
$ws.Range("E28").Value = "  +3.34%  "

$ws.Range("D29").Value = "1.385"
$ws.Range("E29").Value = "  +2.58%  "

$ws.Range("D30").Value = "0.09902"
$ws.Range("E30").Value = "  +2.47%  "

$ws.Range("D31").Value = "4.389"
$ws.Range("E31").Value = "  -0.60%  "

$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").Value = "4.069"
$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("D34").Value = "0.04697"
$ws.Range("E34").Value = "  -0.17%  "

$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  +2.81%  "

$ws.Range("D36").Value = "0.7039"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "2.707"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").Value = "0.01886"
$ws.Range("E38").Value = "  +1.40%  "

$ws.Range("D39").Value = "2.635"
$ws.Range("E39").Value = "  +4.08%  "

$ws.Range("D40").Value = "6.261"
$ws.Range("E40").Value = "  -1.73%  "

$ws.Range("D41").Value = "74.34"
$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").Value = "0.8486"

$ws.Range("D44").Value = "0.4168"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").Value = "0.9994"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").Value = "103.57"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").Value = "965.33"
$ws.Range("E47").Value = "  +3.04%  "

$ws.Range("D48").Value = "7.155"
$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("D49").Value = "9.238"
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("D51").Value = "0.05657"
$ws.Range("E51").Value = "  +0.40%  "

$priceRange.Style = "Normal"
